$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.785.50"
$ws.Range("E2").Value = "  -1.06%  "
$ws.Range("D3").Value = "2.090.67"
$ws.Range("E3").Value = "  +1.91%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'245.09"
$ws.Range("E5").Value = "  -1.40%  "
$ws.Range("D6").Value = "'0.653"
$ws.Range("E6").Value = "  -1.76%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "'54.33"
$ws.Range("E8").Value = "  -5.54%  "
$ws.Range("D9").Value = "'58.78"
$ws.Range("E9").Value = "  -2.10%  "
$ws.Range("E10").Value = "  -4.04%  "
$ws.Range("E11").Value = "  -1.80%  "
$ws.Range("E12").Value = "  +1.00%  "
$ws.Range("D13").Value = "'0.909"
$ws.Range("E13").Value = "  +4.57%  "
$ws.Range("D14").Value = "'15.05"
$ws.Range("E14").Value = "  -5.18%  "
$ws.Range("D15").Value = "2.392.98"
$ws.Range("E15").Value = "  +1.99%  "
$ws.Range("D16").Value = "'5.51"
$ws.Range("E16").Value = "  -2.95%  "
$ws.Range("D17").Value = "2.120.32"
$ws.Range("E17").Value = "  +3.45%  "
$ws.Range("D18").Value = "36.730.98"
$ws.Range("E18").Value = "  -1.26%  "
$ws.Range("E19").Value = "  -4.87%  "
$ws.Range("D20").Value = "'72.71"
$ws.Range("E20").Value = "  -2.57%  "
$ws.Range("D21").Value = "0.0₃0884"
$ws.Range("E21").Value = "  -0.80%  "
$ws.Range("D22").Value = "'5.45"
$ws.Range("E22").Value = "  +1.44%  "
$ws.Range("D23").Value = "'239.04"
$ws.Range("E23").Value = "  +0.85%  "
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").Value = "'2.39"
$ws.Range("E25").Value = "  -3.11%  "
$ws.Range("D26").Value = "'9.79"
$ws.Range("E26").Value = "  +3.29%  "
$ws.Range("E27").Value = "  -0.98%  "
$ws.Range("D28").Value = "'167.45"
$ws.Range("E28").Value = "  -1.04%  "
$ws.Range("E29").Value = "  +2.88%  "
$ws.Range("D30").Value = "'5.32"
$ws.Range("E30").Value = "  +10.43%  "
$ws.Range("E31").Value = "  -1.11%  "
$ws.Range("E32").Value = "  +4.93%  "
$ws.Range("D33").Value = "'4.71"
$ws.Range("E33").Value = "  +5.17%  "
$ws.Range("E34").Value = "  -1.37%  "
$ws.Range("D35").Value = "'2.41"
$ws.Range("E35").Value = "  +7.21%  "
$ws.Range("E36").Value = "  +0.21%  "
$ws.Range("E37").Value = "  +3.86%  "
$ws.Range("D38").Value = "'0.0829"
$ws.Range("E38").Value = "  -6.96%  "
$ws.Range("E39").Value = "  -5.07%  "
$ws.Range("E40").Value = "  +1.71%  "
$ws.Range("E41").Value = "  -0.87%  "
$ws.Range("B42").Value = "THORChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D42").Value = "'4.87"
$ws.Range("E42").Value = "  -7.77%  "
$ws.Range("B43").Value = "Cronos"
$ws.Range("C43").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D43").Value = "'0.0954"
$ws.Range("E43").Value = "  -3.62%  "
$ws.Range("D44").Value = "'96.38"
$ws.Range("E44").Value = "  +0.57%  "
$ws.Range("D45").Value = "'2.86"
$ws.Range("E45").Value = "  -9.44%  "
$ws.Range("D46").Value = "'16.03"
$ws.Range("E46").Value = "  -7.03%  "
$ws.Range("D47").Value = "1.379.34"
$ws.Range("E47").Value = "  +8.82%  "
$ws.Range("D48").Value = "'7.34"
$ws.Range("E48").Value = "  +7.75%  "
$ws.Range("D49").Value = "'2.44"
$ws.Range("E49").Value = "  +0.61%  "
$ws.Range("E50").Value = "  +1.63%  "
$ws.Range("D51").Value = "2.276.06"
$ws.Range("E51").Value = "  +2.16%  "
